# Scheduled-runner update: refresh market-board derived columns (H-N)
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR Leve-profit sheets.
#
# $wb / $excel are pre-seeded by the harness; the workbook is already open.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 751.42426
$ws.Range("I33").Value = 514.7037
$ws.Range("K33").Value = 514.7037
$ws.Range("M33").Value = -285.7037

$ws.Range("H112").Value = 5956.2407
$ws.Range("J112").Value = 6461.9795
$ws.Range("L112").Value = 19385.9385
$ws.Range("N112").Value = -21601.9385

$ws.Range("H121").Value = 2800
$ws.Range("I121").Value = 2200
$ws.Range("J121").Value = 3000
$ws.Range("K121").Value = 6600
$ws.Range("L121").Value = 9000
$ws.Range("M121").Value = -4853
$ws.Range("N121").Value = -12494

$ws.Range("H125").Value = 2156.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2156.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 19408.5
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -24328.5

$ws.Range("H138").Value = 1795.84
$ws.Range("I138").Value = 626.82355
$ws.Range("J138").Value = 2035.2771
$ws.Range("K138").Value = 1880.47065
$ws.Range("L138").Value = 6105.8313
$ws.Range("M138").Value = 3259.52935
$ws.Range("N138").Value = -16385.8313

$ws.Range("H141").Value = 5435.525
$ws.Range("I141").Value = 1972.25
$ws.Range("J141").Value = 10630.4375
$ws.Range("K141").Value = 5916.75
$ws.Range("L141").Value = 31891.3125
$ws.Range("M141").Value = -736.75
$ws.Range("N141").Value = -42251.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 961.0909
$ws.Range("I2").Value = 1027.3529
$ws.Range("J2").Value = 735.8
$ws.Range("K2").Value = 1027.3529
$ws.Range("L2").Value = 735.8
$ws.Range("M2").Value = -914.3529000000001
$ws.Range("N2").Value = -961.8

$ws.Range("H110").Value = 2527.75
$ws.Range("I110").Value = 2527.75
$ws.Range("K110").Value = 2527.75
$ws.Range("M110").Value = -482.75

$ws.Range("H116").Value = 961.0909
$ws.Range("I116").Value = 1027.3529
$ws.Range("J116").Value = 735.8
$ws.Range("K116").Value = 1027.3529
$ws.Range("L116").Value = 735.8
$ws.Range("M116").Value = 1266.6471
$ws.Range("N116").Value = -5323.8

$ws.Range("H123").Value = 40419
$ws.Range("J123").Value = 40419
$ws.Range("L123").Value = 40419
$ws.Range("N123").Value = -50219

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 961.0909
$ws.Range("I3").Value = 1027.3529
$ws.Range("J3").Value = 735.8
$ws.Range("K3").Value = 1027.3529
$ws.Range("L3").Value = 735.8
$ws.Range("M3").Value = -913.3529000000001
$ws.Range("N3").Value = -963.8

$ws.Range("H105").Value = 7814969
$ws.Range("I105").Value = 10418858
$ws.Range("J105").Value = 3300
$ws.Range("K105").Value = 10418858
$ws.Range("L105").Value = 3300
$ws.Range("M105").Value = -10417111
$ws.Range("N105").Value = -6794

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 84615464
$ws.Range("I14").Value = 84615464
$ws.Range("K14").Value = 253846392
$ws.Range("M14").Value = -253846219

$ws.Range("H113").Value = 822.17145
$ws.Range("I113").Value = 444
$ws.Range("K113").Value = 1332
$ws.Range("M113").Value = 838

$ws.Range("H131").Value = 1140.9412
$ws.Range("J131").Value = 1140.9412
$ws.Range("L131").Value = 3422.8236
$ws.Range("N131").Value = -13502.8236

$ws.Range("H139").Value = 2816.3333
$ws.Range("I139").Value = 3133.3333
$ws.Range("J139").Value = 2657.8333
$ws.Range("K139").Value = 9399.999899999999
$ws.Range("L139").Value = 7973.499899999999
$ws.Range("M139").Value = -4259.999899999999
$ws.Range("N139").Value = -18253.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 31000
$ws.Range("J68").Value = 31000
$ws.Range("L68").Value = 31000
$ws.Range("N68").Value = -32622

$ws.Range("H71").Value = 31000
$ws.Range("J71").Value = 31000
$ws.Range("L71").Value = 93000
$ws.Range("N71").Value = -101112

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3500.8
$ws.Range("I40").Value = 3234.6667
$ws.Range("J40").Value = 3900
$ws.Range("K40").Value = 3234.6667
$ws.Range("L40").Value = 3900
$ws.Range("M40").Value = -3098.6667
$ws.Range("N40").Value = -4172

$ws.Range("H55").Value = 835.9259
$ws.Range("I55").Value = 437.875
$ws.Range("J55").Value = 1003.5263
$ws.Range("K55").Value = 437.875
$ws.Range("L55").Value = 1003.5263
$ws.Range("M55").Value = -264.875
$ws.Range("N55").Value = -1349.5263

$ws.Range("H61").Value = 4878.5
$ws.Range("I61").Value = 4742.2
$ws.Range("K61").Value = 4742.2
$ws.Range("M61").Value = -4540.2

$ws.Range("H64").Value = 90000
$ws.Range("J64").Value = 90000
$ws.Range("L64").Value = 90000
$ws.Range("N64").Value = -90450

$ws.Range("H67").Value = 90000
$ws.Range("J67").Value = 90000
$ws.Range("L67").Value = 90000
$ws.Range("N67").Value = -91560

$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100676

$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102340

$ws.Range("H113").Value = 4878.5
$ws.Range("I113").Value = 4742.2
$ws.Range("K113").Value = 4742.2
$ws.Range("M113").Value = -2572.2

$ws.Range("H140").Value = 50214
$ws.Range("J140").Value = 50214
$ws.Range("L140").Value = 50214
$ws.Range("N140").Value = -60574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496

$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716

$ws.Range("H70").Value = 79105
$ws.Range("J70").Value = 79105
$ws.Range("L70").Value = 79105
$ws.Range("N70").Value = -79735

$ws.Range("H73").Value = 79105
$ws.Range("J73").Value = 79105
$ws.Range("L73").Value = 79105
$ws.Range("N73").Value = -81289

$ws.Range("H81").Value = 5311.2666
$ws.Range("I81").Value = 5305.8335
$ws.Range("J81").Value = 5333
$ws.Range("K81").Value = 10611.667
$ws.Range("L81").Value = 10666
$ws.Range("M81").Value = -9550.666999999999
$ws.Range("N81").Value = -12788

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""

$ws.Range("H84").Value = 5311.2666
$ws.Range("I84").Value = 5305.8335
$ws.Range("J84").Value = 5333
$ws.Range("K84").Value = 53058.335
$ws.Range("L84").Value = 53330
$ws.Range("M84").Value = -47754.335
$ws.Range("N84").Value = -63938

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""

$ws.Range("H123").Value = 12396.5
$ws.Range("I123").Value = 10390
$ws.Range("J123").Value = 22429
$ws.Range("K123").Value = 10390
$ws.Range("L123").Value = 22429
$ws.Range("M123").Value = -5490
$ws.Range("N123").Value = -32229

$ws.Range("H126").Value = 1447.3636
$ws.Range("I126").Value = 1197.0769
$ws.Range("J126").Value = 1808.8889
$ws.Range("K126").Value = 3591.2307
$ws.Range("L126").Value = 5426.6667
$ws.Range("M126").Value = -1121.2307
$ws.Range("N126").Value = -10366.6667
